# Auto-generated Excel COM-interop script
# Applies per-cell numeric updates to the 'Aegis_Profits' leve-profit tables
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR, as produced by the
# scheduled market-data refresh run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62: The Mustache Suits Him | Enchanted Mythrite Ink
$ws.Range("H62").Value = 2600.4167
$ws.Range("J62").Value = 2635.3333
$ws.Range("L62").Value = 2635.3333
$ws.Range("N62").Value = -3883.3333

# Row 65: Forgery of Convenience (L) | Enchanted Mythrite Ink
$ws.Range("H65").Value = 2600.4167
$ws.Range("J65").Value = 2635.3333
$ws.Range("L65").Value = 13176.6665
$ws.Range("N65").Value = -19416.6665

# Row 132: Fast-forwarding Flora | Growth Formula Lambda
$ws.Range("H132").Value = 5687376.5
$ws.Range("I132").Value = 6255639
$ws.Range("J132").Value = 4749.75
$ws.Range("K132").Value = 18766917
$ws.Range("L132").Value = 14249.25
$ws.Range("M132").Value = -18764387
$ws.Range("N132").Value = -19309.25

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots | Bronze Ingot
$ws.Range("H2").Value = 35601.69
$ws.Range("I2").Value = 1179.12
$ws.Range("J2").Value = 250742.75
$ws.Range("K2").Value = 1179.12
$ws.Range("L2").Value = 250742.75
$ws.Range("M2").Value = -1066.12
$ws.Range("N2").Value = -250968.75

# Row 114: A New Regular | Bluespirit Gauntlets of Fending
$ws.Range("H114").Value = 28000
$ws.Range("J114").Value = 28000
$ws.Range("L114").Value = 28000
$ws.Range("N114").Value = -36678

# Row 116: No Scope | Titanbronze Ingot
$ws.Range("H116").Value = 35601.69
$ws.Range("I116").Value = 1179.12
$ws.Range("J116").Value = 250742.75
$ws.Range("K116").Value = 1179.12
$ws.Range("L116").Value = 250742.75
$ws.Range("M116").Value = 1114.88
$ws.Range("N116").Value = -255330.75

# Row 122: Haste for High Durium | High Durium Nugget
$ws.Range("H122").Value = 2014.6818
$ws.Range("I122").Value = 1754.8948
$ws.Range("K122").Value = 5264.6844
$ws.Range("M122").Value = -2814.6844

# Row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws.Range("H132").Value = 22603.037
$ws.Range("I132").Value = 26607.227
$ws.Range("J132").Value = 4984.6
$ws.Range("K132").Value = 79821.681
$ws.Range("L132").Value = 14953.8
$ws.Range("M132").Value = -77291.681
$ws.Range("N132").Value = -20013.8

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells | Bronze Ingot
$ws.Range("H3").Value = 35601.69
$ws.Range("I3").Value = 1179.12
$ws.Range("J3").Value = 250742.75
$ws.Range("K3").Value = 1179.12
$ws.Range("L3").Value = 250742.75
$ws.Range("M3").Value = -1065.12
$ws.Range("N3").Value = -250970.75

# Row 26: Unseamly Conditions | Iron Pickaxe
$ws.Range("H26").Value = 7256.2
$ws.Range("I26").Value = 7256.2
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 7256.2
$ws.Range("L26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -6964.2

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value = 35374.395
$ws.Range("I31").Value = 759.0833
$ws.Range("J31").Value = 48773.87
$ws.Range("K31").Value = 759.0833
$ws.Range("L31").Value = 48773.87
$ws.Range("M31").Value = -464.0833
$ws.Range("N31").Value = -49363.87

# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value = 35374.395
$ws.Range("I34").Value = 759.0833
$ws.Range("J34").Value = 48773.87
$ws.Range("K34").Value = 759.0833
$ws.Range("L34").Value = 48773.87
$ws.Range("M34").Value = -557.0833
$ws.Range("N34").Value = -49177.87

$ws = $wb.Worksheets.Item("CUL")
# Row 107: Slippery Service | Frantoio Oil
$ws.Range("H107").Value = 787.0263
$ws.Range("J107").Value = 815.13336
$ws.Range("L107").Value = 2445.40008
$ws.Range("N107").Value = -6285.40008

# Row 109: Cure for What Ails | Purple Carrot Juice
$ws.Range("H109").Value = 3151.973
$ws.Range("I109").Value = 2364.111
$ws.Range("J109").Value = 3405.2144
$ws.Range("K109").Value = 7092.333
$ws.Range("L109").Value = 10215.6432
$ws.Range("M109").Value = -6052.333
$ws.Range("N109").Value = -12295.6432

# Row 110: His Dark Utensils | Spaghetti al Nero
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").ClearContents()

# Row 119: Super Dark Times | Risotto al Nero
$ws.Range("H119").Value = 500300
$ws.Range("I119").Value = 500300
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 1500900
$ws.Range("L119").Value = 0
$ws.Range("M119").ClearContents()
$ws.Range("N119").Value = -1496062

# Row 121: A Cookie for Your Troubles | Coffee Biscuit
$ws.Range("H121").Value = 4373.3555
$ws.Range("I121").Value = 7923.1665
$ws.Range("J121").Value = 3827.2307
$ws.Range("K121").Value = 23769.4995
$ws.Range("L121").Value = 11481.6921
$ws.Range("M121").Value = -22459.4995
$ws.Range("N121").Value = -14101.6921

# Row 129: Comfort Food | Yakow Moussaka
$ws.Range("H129").Value = 178269.94
$ws.Range("I129").Value = 9906
$ws.Range("J129").Value = 222190.95
$ws.Range("K129").Value = 29718
$ws.Range("L129").Value = 666572.8500000001
$ws.Range("M129").Value = -24718
$ws.Range("N129").Value = -676572.8500000001

# Row 130: Blast from the Pasta | The Noodles of Elpis
$ws.Range("H130").Value = 1666.6666
$ws.Range("I130").Value = 1000
$ws.Range("J130").Value = 3000
$ws.Range("K130").Value = 3000
$ws.Range("L130").Value = 9000
$ws.Range("M130").Value = 2020
$ws.Range("N130").Value = -19040

# Row 131: The Mountain Steeped | Tsai tou Vounou
$ws.Range("H131").Value = 813.9
$ws.Range("J131").Value = 857.47253
$ws.Range("L131").Value = 2572.41759
$ws.Range("N131").Value = -12652.41759

$ws = $wb.Worksheets.Item("CRP")
# Row 6: Got Your Back | Square Maple Shield
$ws.Range("H6").Value = 5999.5
$ws.Range("I6").Value = 3499
$ws.Range("J6").Value = 8500
$ws.Range("K6").Value = 3499
$ws.Range("L6").Value = 8500
$ws.Range("M6").Value = -3386
$ws.Range("N6").Value = -8726

$ws = $wb.Worksheets.Item("BSM")
# Row 16: Port of Call: Ul'dah | Bronze Knuckles
$ws.Range("H16").Value = 5999.5
$ws.Range("I16").Value = 3499
$ws.Range("J16").Value = 8500
$ws.Range("K16").Value = 3499
$ws.Range("L16").Value = 8500
$ws.Range("M16").Value = -3249
$ws.Range("N16").Value = -9000

# Row 19: Twice as Slice | Spiked Bronze Labrys
$ws.Range("H19").Value = 500
$ws.Range("I19").Value = 500
$ws.Range("K19").Value = 500
$ws.Range("M19").Value = -212

$ws = $wb.Worksheets.Item("GSM")
# Row 21: Forever 21K | Brass Ring
$ws.Range("H21").Value = 12000.777
$ws.Range("J21").Value = 12000.777
$ws.Range("L21").Value = 12000.777
$ws.Range("N21").Value = -12346.777

# Row 22: Bad to the Bone | Brass Circlet (Sunstone)
$ws.Range("H22").Value = 1509
$ws.Range("J22").Value = 1509
$ws.Range("L22").Value = 1509
$ws.Range("N22").Value = -2567

# Row 30: Dog Tags Are for Dogs | Brass Ring
$ws.Range("H30").Value = 12000.777
$ws.Range("J30").Value = 12000.777
$ws.Range("L30").Value = 12000.777
$ws.Range("N30").Value = -12210.777

# Row 46: Burning the Midnight Oil | Fire Brand
$ws.Range("H46").Value = 12433.111
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 12433.111
$ws.Range("K46").Value = 0
$ws.Range("L46").ClearContents()
$ws.Range("M46").Value = 12433.111
$ws.Range("N46").Value = -12745.111

# Row 97: If I'd a Koppranickel for Every Time... | Koppranickel Ingot
$ws.Range("H97").Value = 34484710
$ws.Range("I97").Value = 50002140
$ws.Range("J97").Value = 1535.7778
$ws.Range("K97").Value = 50002140
$ws.Range("L97").Value = 1535.7778
$ws.Range("M97").Value = -50001644
$ws.Range("N97").Value = -2527.7778

# Row 103: Ring in the New | Azurite Ring of Fending
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").ClearContents()
$ws.Range("N103").Value = 0

# Row 122: Awarding Academic Excellence | Ametrine
$ws.Range("H122").Value = 2743.0908
$ws.Range("I122").Value = 2154.3333
$ws.Range("J122").Value = 5392.5
$ws.Range("K122").Value = 6462.999899999999
$ws.Range("L122").Value = 16177.5
$ws.Range("M122").Value = -4012.999899999999
$ws.Range("N122").Value = -21077.5

$ws = $wb.Worksheets.Item("LTW")
# Row 76: Dragoon Drop Rate | Dhalmelskin Breeches of Maiming
$ws.Range("H76").Value = 9166.333000000001
$ws.Range("J76").Value = 9166.333000000001
$ws.Range("L76").Value = 9166.333000000001
$ws.Range("N76").Value = -9842.333000000001

# Row 79: Exploiting the Adroit (L) | Dhalmelskin Breeches of Maiming
$ws.Range("H79").Value = 9166.333000000001
$ws.Range("J79").Value = 9166.333000000001
$ws.Range("L79").Value = 9166.333000000001
$ws.Range("N79").Value = -11506.333

# Row 106: If the Shoe Fits | Gazelleskin Boots of Casting
$ws.Range("H106").Value = 34076.668
$ws.Range("J106").Value = 34076.668
$ws.Range("L106").Value = 34076.668
$ws.Range("N106").Value = -36600.668

$ws = $wb.Worksheets.Item("WVR")
# Row 113: A Tender Table | Pixie Floss
$ws.Range("H113").Value = 762.53845
$ws.Range("I113").Value = 370.66666
$ws.Range("J113").Value = 880.1
$ws.Range("K113").Value = 1111.99998
$ws.Range("L113").Value = 2640.3
$ws.Range("M113").Value = 1058.00002
$ws.Range("N113").Value = -6980.3

# Row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Range("H132").Value = 4877.1055
$ws.Range("I132").Value = 5338.357
$ws.Range("J132").Value = 3585.6
$ws.Range("K132").Value = 16015.071
$ws.Range("L132").Value = 10756.8
$ws.Range("M132").Value = -13485.071
$ws.Range("N132").Value = -15816.8

